$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared-string text edits (rich text runs; Volume number and report dates) ---
$ws.Range("A8").Value = "Volume 30   Number  48"
$ws.Range("C9").Value = "Report Covering the Week  11/27/2023  Through  12/3/2023"

# --- Column H width change ---
$ws.Columns.Item(8).ColumnWidth = 8.699091

# --- Pure value updates (style/type unchanged) ---
$ws.Range("N16").Value = -89.743589743589
$ws.Range("F19").Value = 3
$ws.Range("G19").Value = 1
$ws.Range("H19").Value = 200
$ws.Range("I19").Value = 48
$ws.Range("K19").Value = 92
$ws.Range("L19").Value = 92
$ws.Range("M19").Value = -21.311475409836
$ws.Range("N19").Value = -72.254335260115
$ws.Range("C21").Value = 2
$ws.Range("F21").Value = 5
$ws.Range("G21").Value = 1
$ws.Range("H21").Value = 400
$ws.Range("I21").Value = 83
$ws.Range("K21").Value = 15.277777777777
$ws.Range("L21").Value = 48.214285714285
$ws.Range("M21").Value = -15.306122448979
$ws.Range("N21").Value = -81.95652173913
$ws.Range("F24").Value = 4
$ws.Range("H24").Value = 300
$ws.Range("J24").Value = 30
$ws.Range("K24").Value = 30
$ws.Range("L24").Value = 18.181818181818
$ws.Range("M24").Value = -61
$ws.Range("J25").Value = 30
$ws.Range("K25").Value = 56.666666666666
$ws.Range("L25").Value = 88
$ws.Range("I30").Value = 4
$ws.Range("L30").Value = 100

# --- Cells that change type/style: number -> shared-string text (reuse style 14 via donor C14/E14) ---
$ws.Range("C14").Copy()
$ws.Range("G15").PasteSpecial(-4163)
$ws.Range("G15").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("C16").PasteSpecial(-4163)
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("G16").PasteSpecial(-4163)
$ws.Range("G16").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("C24").PasteSpecial(-4163)
$ws.Range("C24").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("G26").PasteSpecial(-4163)
$ws.Range("G26").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("F27").PasteSpecial(-4163)
$ws.Range("F27").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("H15").PasteSpecial(-4163)
$ws.Range("H15").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("H16").PasteSpecial(-4163)
$ws.Range("H16").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("H26").PasteSpecial(-4163)
$ws.Range("H26").PasteSpecial(-4122)

# --- Cells that change type/style: shared-string text -> number (reuse style 15/16 via donor J14/K14) ---
$ws.Range("J14").Copy()
$ws.Range("C19").PasteSpecial(-4122)
$ws.Range("C19").Value = 2
$ws.Range("J14").Copy()
$ws.Range("D24").PasteSpecial(-4122)
$ws.Range("D24").Value = 1
$ws.Range("J14").Copy()
$ws.Range("D25").PasteSpecial(-4122)
$ws.Range("D25").Value = 1
$ws.Range("K14").Copy()
$ws.Range("E24").PasteSpecial(-4122)
$ws.Range("E24").Value = -100
$ws.Range("K14").Copy()
$ws.Range("E25").PasteSpecial(-4122)
$ws.Range("E25").Value = -100

$excel.CutCopyMode = 0
